# edit.ps1 -- reproduce the "edit ps5" change to inversesmodn.pptx
#
# Slide 8 ("summary: k is cancellable (mod n) iff k has an inverse (mod n)
# iff gcd(k,n)=1 ... k is relatively prime to n") gets three runs re-colored
# from italic black to non-italic purple (srgbClr 8E00C8):
#   - "cancellable"       (shape "Rectangle 2", 2nd paragraph)
#   - "inverse"           (shape "Rectangle 2", 3rd paragraph)
#   - "relatively prime"  (shape "TextBox 8")
#
# The color is applied via Font.Color.RGB (produces <a:solidFill><a:srgbClr/>
# </a:solidFill>, exactly like the target diff). The italic flag on these
# runs must disappear entirely (not just flip to off): Font.Italic = $false
# only ever rewrites the boolean attribute in place (i="1" -> i="0"), it
# never removes it. To actually drop the attribute, delete the old (italic)
# run and re-insert the same text right after the end of the *preceding*
# (non-italic) run -- the new text merges into that neighboring run's
# formatting (no "i" attribute at all), and re-touching that same span with
# Font.Color afterwards splits out a fresh <a:rPr> (still without "i") that
# carries only the new solid fill -- matching the target XML exactly.

$purple = 0x8E + (0x00 * 256) + (0xC8 * 65536)   # RGB(0x8E,0x00,0xC8) -> 13107342

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(8)

# --- shape "Rectangle 2": "summary: k is cancellable ... k has an inverse ..." ---
$summaryShape = $slide.Shapes.Item(1)
$tr = $summaryShape.TextFrame.TextRange

# "cancellable" -> chars 15-25 (len 11), preceded by " is " at chars 11-14 (len 4)
$old1 = $tr.Characters(15, 11)
$old1.Delete()
$prev1 = $tr.Characters(11, 4)
$prev1.InsertAfter("cancellable") | Out-Null
$new1 = $tr.Characters(15, 11)
$new1.Font.Color.RGB = $purple

# "inverse" -> chars 51-57 (len 7), preceded by " has an " at chars 43-50 (len 8)
$old2 = $tr.Characters(51, 7)
$old2.Delete()
$prev2 = $tr.Characters(43, 8)
$prev2.InsertAfter("inverse") | Out-Null
$new2 = $tr.Characters(51, 7)
$new2.Font.Color.RGB = $purple

# --- shape "TextBox 8": "k is relatively prime to n" ---
$primeShape = $slide.Shapes.Item(4)
$tr2 = $primeShape.TextFrame.TextRange

# "relatively prime" -> chars 6-21 (len 16), preceded by " is " at chars 2-5 (len 4)
$old3 = $tr2.Characters(6, 16)
$old3.Delete()
$prev3 = $tr2.Characters(2, 4)
$prev3.InsertAfter("relatively prime") | Out-Null
$new3 = $tr2.Characters(6, 16)
$new3.Font.Color.RGB = $purple
